$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename return-name values for the vector rows
$ws.Range("B11").Value = "bytevector"
$ws.Range("B12").Value = "charvector"

# Add new "Parameters" column (C) entries for rows that did not have one yet
$ws.Range("C8").Value = "get_glyph_from_byteindex"
$ws.Range("C9").Value = "get_glyphoption_from_byteindex"
$ws.Range("C10").Value = "get_byteslice_from_byteindex"
$ws.Range("C11").Value = "get_bytevector_from_byteindex"
$ws.Range("C12").Value = "get_charvector_from_byteindex"
$ws.Range("C13").Value = "get_glyphvector_from_byteindex"

# Column C needs to widen to fit the new longer strings (~29.71 chars)
$ws.Columns("C").ColumnWidth = 28.8

# Move the active selection to C14, matching the author's final cursor position
$ws.Range("C14").Select()
